$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(121, 8).Value = 4244.5
$ws.Cells.Item(121, 10).Value = 4244.5
$ws.Cells.Item(121, 12).Value = 12733.5
$ws.Cells.Item(121, 14).Value = -16227.5

$ws.Cells.Item(138, 8).Value = 2303.1836
$ws.Cells.Item(138, 9).Value = 2253.5334
$ws.Cells.Item(138, 10).Value = 2381.5789
$ws.Cells.Item(138, 11).Value = 6760.600199999999
$ws.Cells.Item(138, 12).Value = 7144.736699999999
$ws.Cells.Item(138, 13).Value = -1620.600199999999
$ws.Cells.Item(138, 14).Value = -17424.7367

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(23, 8).Value = 29001
$ws.Cells.Item(23, 9).Value = 70006
$ws.Cells.Item(23, 10).Value = 18749.75
$ws.Cells.Item(23, 11).Value = 70006
$ws.Cells.Item(23, 12).Value = 18749.75
$ws.Cells.Item(23, 13).Value = -69747
$ws.Cells.Item(23, 14).Value = -19267.75

$ws.Cells.Item(37, 8).Value = 12018.625
$ws.Cells.Item(37, 9).Value = 10999
$ws.Cells.Item(37, 11).Value = 10999
$ws.Cells.Item(37, 13).Value = -10726

$ws.Cells.Item(53, 8).Value = 16333.333
$ws.Cells.Item(53, 9).Value = 14500
$ws.Cells.Item(53, 11).Value = 14500
$ws.Cells.Item(53, 13).Value = -13818

$ws.Cells.Item(61, 8).Value = 2827.2888
$ws.Cells.Item(61, 9).Value = 2077.7932
$ws.Cells.Item(61, 11).Value = 2077.7932
$ws.Cells.Item(61, 13).Value = -1865.7932

$ws.Cells.Item(63, 8).Value = 9999.5
$ws.Cells.Item(63, 9).Value = 9999.5
$ws.Cells.Item(63, 11).Value = 9999.5
$ws.Cells.Item(63, 13).Value = -9313.5

$ws.Cells.Item(66, 8).Value = 9999.5
$ws.Cells.Item(66, 9).Value = 9999.5
$ws.Cells.Item(66, 11).Value = 49997.5
$ws.Cells.Item(66, 13).Value = -46565.5

$ws.Cells.Item(74, 8).Value = 1266.0571
$ws.Cells.Item(74, 9).Value = 875.9286
$ws.Cells.Item(74, 11).Value = 875.9286
$ws.Cells.Item(74, 13).Value = -1.92859999999996

$ws.Cells.Item(77, 8).Value = 1266.0571
$ws.Cells.Item(77, 9).Value = 875.9286
$ws.Cells.Item(77, 11).Value = 4379.643
$ws.Cells.Item(77, 13).Value = -11.64300000000003

$ws.Cells.Item(102, 8).Value = 2096
$ws.Cells.Item(102, 9).Value = 1895.4286
$ws.Cells.Item(102, 11).Value = 1895.4286
$ws.Cells.Item(102, 13).Value = -273.4286

$ws.Cells.Item(109, 8).Value = 44977.668
$ws.Cells.Item(109, 10).Value = 44977.668
$ws.Cells.Item(109, 12).Value = 44977.668
$ws.Cells.Item(109, 14).Value = -47751.668

$ws.Cells.Item(132, 8).Value = 1923.5625
$ws.Cells.Item(132, 9).Value = 1307.95
$ws.Cells.Item(132, 10).Value = 2949.5833
$ws.Cells.Item(132, 11).Value = 3923.85
$ws.Cells.Item(132, 12).Value = 8848.749899999999
$ws.Cells.Item(132, 13).Value = -1393.85
$ws.Cells.Item(132, 14).Value = -13908.7499

$ws.Cells.Item(136, 8).Value = 2827.2888
$ws.Cells.Item(136, 9).Value = 2077.7932
$ws.Cells.Item(136, 11).Value = 6233.3796
$ws.Cells.Item(136, 13).Value = -3683.3796

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2078.7
$ws.Cells.Item(20, 9).Value = 2041
$ws.Cells.Item(20, 10).Value = 2166.6667
$ws.Cells.Item(20, 11).Value = 2041
$ws.Cells.Item(20, 12).Value = 2166.6667
$ws.Cells.Item(20, 13).Value = -1794
$ws.Cells.Item(20, 14).Value = -2660.6667

$ws.Cells.Item(86, 8).Value = 107539.79
$ws.Cells.Item(86, 9).Value = 2890.4546
$ws.Cells.Item(86, 10).Value = 251432.62
$ws.Cells.Item(86, 11).Value = 2890.4546
$ws.Cells.Item(86, 12).Value = 251432.62
$ws.Cells.Item(86, 13).Value = -1767.4546
$ws.Cells.Item(86, 14).Value = -253678.62

$ws.Cells.Item(89, 8).Value = 107539.79
$ws.Cells.Item(89, 9).Value = 2890.4546
$ws.Cells.Item(89, 10).Value = 251432.62
$ws.Cells.Item(89, 11).Value = 14452.273
$ws.Cells.Item(89, 12).Value = 1257163.1
$ws.Cells.Item(89, 13).Value = -8836.273000000001
$ws.Cells.Item(89, 14).Value = -1268395.1

$ws.Cells.Item(105, 8).Value = 2247.9333
$ws.Cells.Item(105, 9).Value = 2158.5
$ws.Cells.Item(105, 10).Value = 3500
$ws.Cells.Item(105, 11).Value = 2158.5
$ws.Cells.Item(105, 12).Value = 3500
$ws.Cells.Item(105, 13).Value = -411.5
$ws.Cells.Item(105, 14).Value = -6994

$ws.Cells.Item(107, 8).Value = 699.5
$ws.Cells.Item(107, 9).Value = 700
$ws.Cells.Item(107, 10).Value = 699
$ws.Cells.Item(107, 11).Value = 700
$ws.Cells.Item(107, 12).Value = 699
$ws.Cells.Item(107, 13).Value = 1220
$ws.Cells.Item(107, 14).Value = -4539

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 949.8333
$ws.Cells.Item(22, 9).Value = 372
$ws.Cells.Item(22, 10).Value = 1238.75
$ws.Cells.Item(22, 11).Value = 372
$ws.Cells.Item(22, 12).Value = 1238.75
$ws.Cells.Item(22, 13).Value = -22
$ws.Cells.Item(22, 14).Value = -1938.75

$ws.Cells.Item(31, 8).Value = 2539.7917
$ws.Cells.Item(31, 9).Value = 2407.0908
$ws.Cells.Item(31, 10).Value = 2652.077
$ws.Cells.Item(31, 11).Value = 2407.0908
$ws.Cells.Item(31, 12).Value = 2652.077
$ws.Cells.Item(31, 13).Value = -2112.0908
$ws.Cells.Item(31, 14).Value = -3242.077

$ws.Cells.Item(34, 8).Value = 2539.7917
$ws.Cells.Item(34, 9).Value = 2407.0908
$ws.Cells.Item(34, 10).Value = 2652.077
$ws.Cells.Item(34, 11).Value = 2407.0908
$ws.Cells.Item(34, 12).Value = 2652.077
$ws.Cells.Item(34, 13).Value = -2205.0908
$ws.Cells.Item(34, 14).Value = -3056.077

$ws.Cells.Item(58, 8).Value = 1145541.1
$ws.Cells.Item(58, 9).Value = 1500236.8
$ws.Cells.Item(58, 10).Value = 2633.2222
$ws.Cells.Item(58, 11).Value = 1500236.8
$ws.Cells.Item(58, 12).Value = 2633.2222
$ws.Cells.Item(58, 13).Value = -1500033.8
$ws.Cells.Item(58, 14).Value = -3039.2222

$ws.Cells.Item(86, 8).Value = 111113120
$ws.Cells.Item(86, 10).Value = 2947.75
$ws.Cells.Item(86, 12).Value = 2947.75
$ws.Cells.Item(86, 14).Value = -5193.75

$ws.Cells.Item(89, 8).Value = 111113120
$ws.Cells.Item(89, 10).Value = 2947.75
$ws.Cells.Item(89, 12).Value = 14738.75
$ws.Cells.Item(89, 14).Value = -25970.75

$ws.Cells.Item(132, 8).Value = 1937.0476
$ws.Cells.Item(132, 9).Value = 1204.9656
$ws.Cells.Item(132, 11).Value = 3614.8968
$ws.Cells.Item(132, 13).Value = -1084.8968

$ws.Cells.Item(134, 8).Value = 1092.0444
$ws.Cells.Item(134, 9).Value = 1122.35
$ws.Cells.Item(134, 11).Value = 3367.05
$ws.Cells.Item(134, 13).Value = -832.0499999999997

$ws.Cells.Item(136, 8).Value = 1145541.1
$ws.Cells.Item(136, 9).Value = 1500236.8
$ws.Cells.Item(136, 10).Value = 2633.2222
$ws.Cells.Item(136, 11).Value = 4500710.4
$ws.Cells.Item(136, 12).Value = 7899.6666
$ws.Cells.Item(136, 13).Value = -4498160.4
$ws.Cells.Item(136, 14).Value = -12999.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 131.4
$ws.Cells.Item(33, 9).Value = 92.333336
$ws.Cells.Item(33, 10).Value = 190
$ws.Cells.Item(33, 11).Value = 554.000016
$ws.Cells.Item(33, 12).Value = 1140
$ws.Cells.Item(33, 13).Value = -271.000016
$ws.Cells.Item(33, 14).Value = -1706

$ws.Cells.Item(80, 8).Value = 2000
$ws.Cells.Item(80, 10).Value = 2000
$ws.Cells.Item(80, 12).Value = 6000
$ws.Cells.Item(80, 14).Value = -7872

$ws.Cells.Item(83, 8).Value = 2000
$ws.Cells.Item(83, 10).Value = 2000
$ws.Cells.Item(83, 12).Value = 18000
$ws.Cells.Item(83, 14).Value = -27360

$ws.Cells.Item(98, 8).Value = 92.666664
$ws.Cells.Item(98, 9).Value = 92.2
$ws.Cells.Item(98, 10).Value = 95
$ws.Cells.Item(98, 11).Value = 276.6
$ws.Cells.Item(98, 12).Value = 285
$ws.Cells.Item(98, 13).Value = 1221.4
$ws.Cells.Item(98, 14).Value = -3281

$ws.Cells.Item(122, 8).Value = 1037.3334
$ws.Cells.Item(122, 9).Value = 762.25
$ws.Cells.Item(122, 10).Value = 1137.3636
$ws.Cells.Item(122, 11).Value = 6860.25
$ws.Cells.Item(122, 12).Value = 10236.2724
$ws.Cells.Item(122, 13).Value = -4410.25
$ws.Cells.Item(122, 14).Value = -15136.2724

$ws.Cells.Item(131, 8).Value = 773.41
$ws.Cells.Item(131, 10).Value = 804.6923
$ws.Cells.Item(131, 12).Value = 2414.0769
$ws.Cells.Item(131, 14).Value = -12494.0769

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 1997
$ws.Cells.Item(43, 9).Value = 1997
$ws.Cells.Item(43, 11).Value = 1997
$ws.Cells.Item(43, 13).Value = -1846

$ws.Cells.Item(97, 8).Value = 1364.4546
$ws.Cells.Item(97, 9).Value = 1200.9
$ws.Cells.Item(97, 11).Value = 1200.9
$ws.Cells.Item(97, 13).Value = -704.9000000000001

$ws.Cells.Item(102, 8).Value = 3082.077
$ws.Cells.Item(102, 9).Value = 3097.0908
$ws.Cells.Item(102, 10).Value = 2999.5
$ws.Cells.Item(102, 11).Value = 3097.0908
$ws.Cells.Item(102, 12).Value = 2999.5
$ws.Cells.Item(102, 13).Value = -1475.0908
$ws.Cells.Item(102, 14).Value = -6243.5

$ws.Cells.Item(119, 8).Value = 50000
$ws.Cells.Item(119, 10).Value = 50000
$ws.Cells.Item(119, 12).Value = 50000
$ws.Cells.Item(119, 14).Value = -59676

$ws.Cells.Item(122, 8).Value = 1804.5294
$ws.Cells.Item(122, 9).Value = 1456.5834
$ws.Cells.Item(122, 10).Value = 2639.6
$ws.Cells.Item(122, 11).Value = 4369.7502
$ws.Cells.Item(122, 12).Value = 7918.799999999999
$ws.Cells.Item(122, 13).Value = -1919.7502
$ws.Cells.Item(122, 14).Value = -12818.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4058.7144
$ws.Cells.Item(7, 9).Value = 2484.3
$ws.Cells.Item(7, 11).Value = 2484.3
$ws.Cells.Item(7, 13).Value = -2372.3

$ws.Cells.Item(61, 8).Value = 3262.375
$ws.Cells.Item(61, 9).Value = 3149.8333
$ws.Cells.Item(61, 11).Value = 3149.8333
$ws.Cells.Item(61, 13).Value = -2947.8333

$ws.Cells.Item(82, 8).Value = 807.5
$ws.Cells.Item(82, 9).Value = 807.5
$ws.Cells.Item(82, 11).Value = 807.5
$ws.Cells.Item(82, 13).Value = -446.5

$ws.Cells.Item(85, 8).Value = 807.5
$ws.Cells.Item(85, 9).Value = 807.5
$ws.Cells.Item(85, 11).Value = 807.5
$ws.Cells.Item(85, 13).Value = 440.5

$ws.Cells.Item(113, 8).Value = 3262.375
$ws.Cells.Item(113, 9).Value = 3149.8333
$ws.Cells.Item(113, 11).Value = 3149.8333
$ws.Cells.Item(113, 13).Value = -979.8332999999998

$ws.Cells.Item(126, 8).Value = 4058.7144
$ws.Cells.Item(126, 9).Value = 2484.3
$ws.Cells.Item(126, 11).Value = 7452.900000000001
$ws.Cells.Item(126, 13).Value = -4982.900000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1167.3334
$ws.Cells.Item(107, 9).Value = 1000.5
$ws.Cells.Item(107, 11).Value = 3001.5
$ws.Cells.Item(107, 13).Value = -1081.5

$ws.Cells.Item(132, 8).Value = 1555.3334
$ws.Cells.Item(132, 9).Value = 1123.8334
$ws.Cells.Item(132, 10).Value = 3281.3333
$ws.Cells.Item(132, 11).Value = 3371.5002
$ws.Cells.Item(132, 12).Value = 9843.999899999999
$ws.Cells.Item(132, 13).Value = -841.5001999999999
$ws.Cells.Item(132, 14).Value = -14903.9999
